$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row height tweaks (rows 13 / 14) ---
$ws.Rows.Item(13).RowHeight = 39.75
$ws.Rows.Item(14).RowHeight = 54.75

# --- New entries: rows 41 and 42 (previously blank placeholder rows) ---
# Duplicate the formatting of the row above (row 40, the last filled entry)
# into the two rows we are about to populate.
$ws.Range("A40:G40").Copy($ws.Range("A41:G41"))
$ws.Range("A40:G40").Copy($ws.Range("A42:G42"))

# Row 41: entry #10 - "Работа над отчетом"
$ws.Cells.Item(41, 1).Value = 10
$ws.Cells.Item(41, 2).Value = "Работа над отчетом"
$ws.Cells.Item(41, 3).Value = (Get-Date -Year 2024 -Month 7 -Day 5 -Hour 9 -Minute 0 -Second 0)
$ws.Cells.Item(41, 4).Value = (Get-Date -Year 2024 -Month 7 -Day 5 -Hour 17 -Minute 0 -Second 0)
$ws.Cells.Item(41, 5).Formula = "=D41-C41"
$ws.Cells.Item(41, 6).Value = $null
$ws.Cells.Item(41, 7).Formula = "=SUM(E41)"

# Row 42: entry #11 - "Работа над отчетом"
$ws.Cells.Item(42, 1).Value = 11
$ws.Cells.Item(42, 2).Value = "Работа над отчетом"
$ws.Cells.Item(42, 3).Value = (Get-Date -Year 2024 -Month 7 -Day 5 -Hour 9 -Minute 0 -Second 0)
$ws.Cells.Item(42, 4).Value = (Get-Date -Year 2024 -Month 7 -Day 5 -Hour 15 -Minute 0 -Second 0)
$ws.Cells.Item(42, 5).Formula = "=D42-C42"
$ws.Cells.Item(42, 6).Value = $null
$ws.Cells.Item(42, 7).Formula = "=SUM(E42)"

# Row heights for the two new rows
$ws.Rows.Item(41).RowHeight = 36
$ws.Rows.Item(42).RowHeight = 15.75

# --- View / selection state ---
$ws.Range("I40").Select()
